$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells holding numeric-looking values that must remain stored as text
# (matching the original inlineStr formatting with exact digits/trailing zeros).
$textNumericCells = @{
    "D2" = "262.16"
    "G2" = "4"
    "D3" = "22.93"
    "G3" = "4"
    "G4" = "4"
    "D5" = "0.06103"
    "G5" = "4"
    "D6" = "6.735"
    "G6" = "4"
    "D7" = "3.462"
    "G7" = "4"
    "D8" = "1.366"
    "G8" = "4"
    "D9" = "0.7981"
    "G9" = "4"
    "G10" = "4"
    "D11" = "0.08114"
    "G11" = "4"
    "D12" = "0.03467"
    "G12" = "4"
    "D13" = "0.03048"
    "G13" = "4"
    "D14" = "0.09320"
    "G14" = "4"
    "D15" = "3.861"
    "G15" = "4"
    "D16" = "0.001700"
    "G16" = "4"
    "D17" = "0.04784"
    "G17" = "4"
    "D18" = "0.0006148"
    "G18" = "4"
    "D19" = "0.006238"
    "G19" = "4"
    "D20" = "0.001091"
    "G20" = "4"
    "D21" = "0.003809"
    "G21" = "4"
    "G22" = "4"
    "D23" = "3.704"
    "G23" = "4"
    "D24" = "2.205"
    "G24" = "4"
    "G25" = "4"
    "D26" = "0.1252"
    "G26" = "4"
    "D27" = "0.0003201"
    "G27" = "4"
    "G28" = "4"
    "G29" = "4"
    "G30" = "4"
    "G31" = "4"
    "G32" = "4"
    "G33" = "4"
    "G34" = "4"
    "G35" = "4"
    "G36" = "4"
    "G37" = "4"
    "G38" = "4"
    "G39" = "4"
    "D40" = "0.04608"
    "G40" = "4"
    "D41" = "0.1119"
    "G41" = "4"
    "D42" = "0.003130"
    "G42" = "4"
    "D43" = "0.003331"
    "G43" = "4"
    "D44" = "0.01023"
    "G44" = "4"
    "D45" = "0.002969"
    "G45" = "4"
    "D46" = "0.00005926"
    "G46" = "4"
    "G47" = "4"
    "D48" = "0.6997"
    "G48" = "4"
    "D49" = "0.07374"
    "G49" = "4"
    "D50" = "0.00002099"
    "G50" = "4"
    "G51" = "4"
}
foreach ($ref in $textNumericCells.Keys) {
    $cell = $ws.Range($ref)
    $cell.NumberFormat = "@"
    $cell.Value = $textNumericCells[$ref]
}

# Plain text cells (coin names, URLs, volume labels)
$textCells = @{
    "B41" = "BKEXToken"
    "C41" = "https://coinranking.com/coin/IPeThtYgk+bkextoken-bkk"
    "E41" = "40BKEXTokenBKK"
    "B43" = "KickToken"
    "C43" = "https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick"
    "E43" = "42KickTokenKICK"
}
foreach ($ref in $textCells.Keys) {
    $ws.Range($ref).Value = $textCells[$ref]
}
